$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 78: correct the date/time value in column A ---
$ws.Range("A78").Value = 45449.2916666667

# --- Add new row 79 with a new OHLC data point ---

# Column A: date serial value, re-using the date/time style from A78
$ws.Range("A78").Copy($ws.Range("A79"))
$ws.Range("A79").Value = 45450.583125

# Column B: volume
$ws.Range("B79").Value = 600

# Columns C-F: high/low/open/close (all equal for this row)
$ws.Range("C79").Value = 6.28000020980835
$ws.Range("D79").Value = 6.28000020980835
$ws.Range("E79").Value = 6.28000020980835
$ws.Range("F79").Value = 6.28000020980835

# Column G: adj_close stored as text (shared string), matching existing
# convention in the sheet where this column holds text representations
# of the numeric value. A leading apostrophe forces Excel to store the
# value as text instead of auto-converting it to a number; resetting
# the style back to Normal afterwards avoids leaving a quote-prefix
# style applied to the cell (matching the un-styled cells elsewhere in
# this column).
$ws.Range("G79").Value = "'6.28000020980835"
$ws.Range("G79").Style = "Normal"

# Column H: ticker
$ws.Range("H79").Value = "PAL.MI"
